$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "seating flexed no support" row entirely (old row 3). This
# shifts every subsequent row up by one and lets Excel collapse the
# chart's non-contiguous series references back into simple ranges.
$ws.Rows.Item(3).Delete()

# Updated AMS [N] measurements (column C) after the AMMR rework. Row
# numbers below are the *new* (post-deletion) row numbers.
$ws.Range("C2").Value = 52
$ws.Range("C3").Value = 245
$ws.Range("C4").Value = 452
$ws.Range("C5").Value = 506
$ws.Range("C6").Value = 1066
$ws.Range("C7").Value = 1242
$ws.Range("C8").Value = 2640
$ws.Range("C9").Value = 2120

# Move the active selection like the saved workbook shows.
$ws.Range("O15").Select()
